$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A149 literal value
$ws.Range("A149").Value = 309

# Update dependent formulas on D237, D241, D245
$ws.Range("D237").Formula = "=2*B30+B18"
$ws.Range("D241").Formula = "=2*B30+C18"
$ws.Range("D245").Formula = "=2*B30+D18"

# Update sheet view (selection / scroll position)
$excel.ActiveWindow.ScrollRow = 224
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A267").Select()
